$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Replace every "deuteron" value in the target column (G2:G13) with "d"
$target = $ws.Range("G2:G13")
foreach ($cell in $target.Cells) {
    if ($cell.Value2 -eq "deuteron") {
        $cell.Value = "d"
    }
}

# Make the header row bold and centered
$header = $ws.Range("A1:K1")
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108

# Match the saved selection state (F16)
$ws.Range("F16").Select() | Out-Null
